$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Top performers
$ws.Range("B2").Value = 0.3750000000000001

$ws.Range("A3").Value = "VSS1R"
$ws.Range("B3").Value = 0.0984848484848484

$ws.Range("A4").Value = "SKN1T"
$ws.Range("B4").Value = 0.09090909090909087

$ws.Range("A5").Value = "UTR1L"
$ws.Range("B5").Value = 0.05000000000000004

$ws.Range("A6").Value = "PKG1T"
$ws.Range("B6").Value = 0.04444444444444431

# Worst performers
$ws.Range("A8").Value = "TPD1T"
$ws.Range("B8").Value = -0.3870967741935484

$ws.Range("A9").Value = "BLT1T"
$ws.Range("B9").Value = -0.2192691029900332

$ws.Range("A10").Value = "INC1L"
$ws.Range("B10").Value = -0.04878048780487796

$ws.Range("A11").Value = "OLF1R"
$ws.Range("B11").Value = -0.03380281690140836

$ws.Range("A12").Value = "NTU1L"
$ws.Range("B12").Value = -0.03267973856209151

# Dates (kept as plain text, not Excel date serials)
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "2019.07.26"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "2019.07.19"

# Ratio label
$ws.Range("A21").Value = "ROE"

# Table of companies (U/V columns)
$ws.Range("U24").Value = "PATA Saldus"
$ws.Range("V24").Value = 0.4315683497470772

$ws.Range("U25").Value = "Silvano Fashion Group"
$ws.Range("V25").Value = 0.30216277307756

$ws.Range("U26").Value = "Tallinna Vesi"
$ws.Range("V26").Value = 0.2561233634353408

$ws.Range("V27").Value = 0.2324629178656131

$ws.Range("U28").Value = "MADARA Cosmetics"
$ws.Range("V28").Value = 0.1933831470361544
